# Refresh the crypto price/volume snapshot (Price column D, Volume(1h) column E).
# All D/E cells in the sheet are stored as text (t="inlineStr" in the source
# workbook), so for Price values that parse as a plain number (e.g. "0.997",
# "214.60") we lead with an apostrophe to force a text entry - otherwise
# Excel would silently coerce the string to a numeric cell and drop
# formatting such as trailing zeros (e.g. "214.60" -> 214.6). Values that
# contain extra punctuation (e.g. "26.575.59") already fail numeric parsing
# and don't need the prefix; Volume(1h) values always carry a "%" and
# padding spaces, so they're never at risk either.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.575.59'
$ws.Range("E2").Value = '  +1.23%  '
$ws.Range("D3").Value = '1.622.16'
$ws.Range("E3").Value = '  +1.98%  '
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").Value = "'214.60"
$ws.Range("E5").Value = '  +1.06%  '
$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = '  +0.80%  '
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("E8").Value = '  +0.42%  '
$ws.Range("D9").Value = "'0.0612"
$ws.Range("E9").Value = '  +0.61%  '
$ws.Range("D10").Value = "'19.38"
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("D11").Value = "'0.0856"
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("D12").Value = '1.847.56'
$ws.Range("E12").Value = '  +1.88%  '
$ws.Range("D13").Value = '1.623.60'
$ws.Range("E13").Value = '  +1.85%  '
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").Value = "'0.514"
$ws.Range("E15").Value = '  -1.31%  '
$ws.Range("D16").Value = "'64.88"
$ws.Range("E16").Value = '  +0.71%  '
$ws.Range("D17").Value = '26.546.61'
$ws.Range("E17").Value = '  +1.16%  '
$ws.Range("D18").Value = "'231.85"
$ws.Range("E18").Value = '  +8.36%  '
$ws.Range("D19").Value = '0.0₃0729'
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("D20").Value = "'7.60"
$ws.Range("E20").Value = '  +2.17%  '
$ws.Range("D21").Value = "'0.997"
$ws.Range("E21").Value = '  -0.38%  '
$ws.Range("E22").Value = '  +2.19%  '
$ws.Range("D23").Value = "'9.15"
$ws.Range("E23").Value = '  +1.60%  '
$ws.Range("D24").Value = "'2.16"
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("E25").Value = '  +1.22%  '
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("D27").Value = "'7.05"
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("E28").Value = '  +2.19%  '
$ws.Range("D29").Value = "'15.69"
$ws.Range("E29").Value = '  +3.18%  '
$ws.Range("D30").Value = "'0.0499"
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").Value = "'1.17"
$ws.Range("E31").Value = '  +0.96%  '
$ws.Range("D32").Value = "'3.24"
$ws.Range("E32").Value = '  +1.45%  '
$ws.Range("D33").Value = '1.444.55'
$ws.Range("E33").Value = '  +8.09%  '
$ws.Range("E34").Value = '  +2.23%  '
$ws.Range("E35").Value = '  -0.91%  '
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("E37").Value = '  -5.19%  '
$ws.Range("D38").Value = "'0.0168"
$ws.Range("E38").Value = '  +0.62%  '
$ws.Range("D39").Value = "'0.840"
$ws.Range("E39").Value = '  +2.80%  '
$ws.Range("E40").Value = '  +2.21%  '
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("E42").Value = '  +2.84%  '
$ws.Range("D43").Value = '1.758.54'
$ws.Range("E43").Value = '  +1.98%  '
$ws.Range("D44").Value = "'0.764"
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("D45").Value = "'62.35"
$ws.Range("E45").Value = '  +0.76%  '
$ws.Range("D46").Value = "'0.920"
$ws.Range("E46").Value = '  -8.80%  '
$ws.Range("E47").Value = '  +3.15%  '
$ws.Range("E48").Value = '  -1.32%  '
$ws.Range("E49").Value = '  +1.21%  '
$ws.Range("E50").Value = '  +0.33%  '
